# Reading_Activity4.docx update
#  1. Add a "_GoBack" bookmark on the (previously) blank paragraph right
#     after the "Program #6" heading.
#  2. Rename the tryAgain variable to guessAgain throughout (3 call sites).
#  3. attemptNumber is now initialised to 0 instead of 1.
#  4. Drop now-unneeded leading spaces on a few statement continuations.
#  5. Change "...prompted to try again?" to "...prompted to guess again?"
#     (split into "gues" / "s" / " again?" runs, matching how Word
#     records an in-place retype). The old "_GoBack" bookmark that used
#     to wrap "try " here is already gone, relocated by step 1.

$d = $word.ActiveDocument

# --- 1. New _GoBack bookmark on the blank paragraph after "Program #6" ---
$blankPara = $d.Paragraphs(148)
[void]$d.Bookmarks.Add("_GoBack", $blankPara.Range)

# --- 2. tryAgain -> guessAgain (all three occurrences) ---
[void]$d.Content.Find.Execute("tryAgain", $false, $false, $false, $false, $false, `
    $true, 1, $false, "guessAgain", 2)

# --- 3. attemptNumber = 1  ->  attemptNumber = 0 -------------------------
$attemptPara = $d.Paragraphs(158)
$r = $attemptPara.Range
$found = $r.Find.Execute(" = 1")
if ($found) {
    $r.Text = " = 0"
}

# --- 4. Trim leading spaces that are no longer needed --------------------
$r = $d.Content
if ($r.Find.Execute(' = "yes"')) {
    $r.Text = '= "yes"'
}

$r = $d.Content
if ($r.Find.Execute(' != "no":')) {
    $r.Text = '!= "no":'
}

$r = $d.Content
if ($r.Find.Execute(" = input('Would you like to guess again? [")) {
    $r.Text = "= input('Would you like to guess again? ["
}

# --- 5. Turn "try again?" into "guess again?" (the old _GoBack bookmark --
#        around "try " was already relocated away from here in step 1,
#        since a document can only have a single "_GoBack" bookmark).
$anchor = $d.Content
[void]$anchor.Find.Execute("prompted to try again?")
$tryStart = $anchor.Start + 12
$tryEnd = $tryStart + 4
$againStart = $tryEnd
$againEnd = $againStart + 6

# "try " -> "gues"  (toggle Bold off/on so the new text is not silently
# re-merged into the identically formatted preceding run)
$rTry = $d.Range($tryStart, $tryEnd)
$rTry.Font.Bold = 1
$rTry.Text = "gues"
$rTry2 = $d.Range($tryStart, $tryStart + 4)
$rTry2.Font.Bold = $false

# "again?" -> "s again?" as its own run
$rAgain = $d.Range($againStart, $againEnd)
$rAgain.Font.Bold = 1
$rAgain.Text = "s again?"
$rAgainFixed = $d.Range($againStart, $againStart + 8)
$rAgainFixed.Font.Bold = $false

# split "s again?" into "s" and " again?"
$rS = $d.Range($againStart, $againStart + 1)
$rS.Font.Bold = 1
$rS.Text = "s"
$rS2 = $d.Range($againStart, $againStart + 1)
$rS2.Font.Bold = $false

Write-Output ("Final paragraph text: " + $d.Paragraphs(181).Range.Text)
